$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds digit-only strings (e.g. "10000"). Excel.Range.Value
# auto-detects those as numbers, but the target stores them as plain
# inline strings, so format just those specific cells as Text first.
$gCells = @("G40","G41","G42","G43","G44","G45","G46","G49")
foreach ($addr in $gCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("B32").Value = "High School"
$ws.Range("C32").Value = "District of Columbia"
$ws.Range("D32").Value = "out-of-state"
$ws.Range("E32").Value = "Education"
$ws.Range("F32").Value = "Less than 1 Year"
$ws.Range("H32").Value = 44027.66921391546
$ws.Range("H32").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B33").Value = "High School"
$ws.Range("C33").Value = "Iowa"
$ws.Range("D33").Value = "out-of-state"
$ws.Range("E33").Value = "Education"
$ws.Range("F33").Value = "Less than 1 Year"
$ws.Range("H33").Value = 44027.66937459796
$ws.Range("H33").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B34").Value = "High School"
$ws.Range("C34").Value = "District of Columbia"
$ws.Range("D34").Value = "out-of-state"
$ws.Range("E34").Value = "Education"
$ws.Range("F34").Value = "Less than 1 Year"
$ws.Range("H34").Value = 44027.66950079981
$ws.Range("H34").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B35").Value = "High School"
$ws.Range("C35").Value = "District of Columbia"
$ws.Range("D35").Value = "out-of-state"
$ws.Range("E35").Value = "Agriculture & Natural Resources"
$ws.Range("F35").Value = "Less than 1 Year"
$ws.Range("H35").Value = 44027.66965874392
$ws.Range("H35").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B36").Value = "High School"
$ws.Range("C36").Value = "District of Columbia"
$ws.Range("D36").Value = "in-state"
$ws.Range("E36").Value = "Agriculture & Natural Resources"
$ws.Range("F36").Value = "Less than 1 Year"
$ws.Range("H36").Value = 44027.66973322234
$ws.Range("H36").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B37").Value = "High School"
$ws.Range("C37").Value = "Iowa"
$ws.Range("D37").Value = "in-state"
$ws.Range("E37").Value = "Humanities & Liberal Arts"
$ws.Range("F37").Value = "Less than 1 Year"
$ws.Range("H37").Value = 44027.66982667609
$ws.Range("H37").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B38").Value = "High School"
$ws.Range("C38").Value = "New York"
$ws.Range("D38").Value = "out-of-state"
$ws.Range("E38").Value = "Engineering"
$ws.Range("F38").Value = "Less than 1 Year"
$ws.Range("H38").Value = 44027.66995940193
$ws.Range("H38").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B39").Value = "High School"
$ws.Range("C39").Value = "New Jersey"
$ws.Range("D39").Value = "in-state"
$ws.Range("E39").Value = "Education"
$ws.Range("F39").Value = "Less than 1 Year"
$ws.Range("H39").Value = 44027.67013829082
$ws.Range("H39").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B40").Value = "College"
$ws.Range("C40").Value = "District of Columbia"
$ws.Range("E40").Value = "Agriculture & Natural Resources"
$ws.Range("G40").Value = "10000"
$ws.Range("H40").Value = 44027.67032781119
$ws.Range("H40").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B41").Value = "College"
$ws.Range("C41").Value = "District of Columbia"
$ws.Range("E41").Value = "Engineering"
$ws.Range("G41").Value = "100000"
$ws.Range("H41").Value = 44027.67051353793
$ws.Range("H41").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B42").Value = "College"
$ws.Range("C42").Value = "New Jersey"
$ws.Range("E42").Value = "Education"
$ws.Range("G42").Value = "25000"
$ws.Range("H42").Value = 44027.6708477892
$ws.Range("H42").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B43").Value = "College"
$ws.Range("C43").Value = "New York"
$ws.Range("E43").Value = "Education"
$ws.Range("G43").Value = "25000"
$ws.Range("H43").Value = 44027.67152561496
$ws.Range("H43").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B44").Value = "College"
$ws.Range("C44").Value = "New York"
$ws.Range("E44").Value = "Engineering"
$ws.Range("G44").Value = "25000"
$ws.Range("H44").Value = 44027.67189825617
$ws.Range("H44").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B45").Value = "College"
$ws.Range("C45").Value = "Michigan"
$ws.Range("E45").Value = "Engineering"
$ws.Range("G45").Value = "100000"
$ws.Range("H45").Value = 44027.67257344551
$ws.Range("H45").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B46").Value = "College"
$ws.Range("C46").Value = "Illinois"
$ws.Range("E46").Value = "Engineering"
$ws.Range("G46").Value = "100000"
$ws.Range("H46").Value = 44027.67270165341
$ws.Range("H46").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B47").Value = "High School"
$ws.Range("C47").Value = "District of Columbia"
$ws.Range("D47").Value = "in-state"
$ws.Range("E47").Value = "Agriculture & Natural Resources"
$ws.Range("F47").Value = "Less than 1 Year"
$ws.Range("H47").Value = 44027.67284026599
$ws.Range("H47").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B48").Value = "High School"
$ws.Range("C48").Value = "New Jersey"
$ws.Range("D48").Value = "in-state"
$ws.Range("E48").Value = "Education"
$ws.Range("F48").Value = "Less than 1 Year"
$ws.Range("H48").Value = 44027.67912193229
$ws.Range("H48").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("B49").Value = "College"
$ws.Range("C49").Value = "Alabama"
$ws.Range("E49").Value = "Agriculture & Natural Resources"
$ws.Range("G49").Value = "25000"
$ws.Range("H49").Value = 44027.68164028682
$ws.Range("H49").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Drop the Text number format back to the default/unstyled format so
# the written G cells carry no explicit style (matches the target,
# which has no s= attribute on those cells), while keeping their string content.
foreach ($addr in $gCells) { $ws.Range($addr).Style = "Normal" }
